$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing) - update B2 and E2 to the new album/song values
$ws.Range("B2").Value = "The hold steady"
$ws.Range("E2").Value = "Boys and girls in america"

# Row 3 (new)
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "sujian stevens"
$ws.Range("C3").Value = "mdex:string"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "oh great white city"

# Row 4 (new)
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "hulk"
$ws.Range("C4").Value = "mdex:string"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "SMASH"

# Update selection to match the target state
$ws.Range("E4").Select()
